# Fix S dock numbering so each route finishes with a docking station.
# Also sync up the summary route dock counts with reality.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- S route (rows 17-28), dock/number counts in column C ---
$ws.Range("C18").Value = 17
$ws.Range("C19").Value = 21
$ws.Range("C20").Value = 26
$ws.Range("C22").Value = 17

# --- S route summary (rows 32-43), dock/number counts in column C ---
$ws.Range("C34").Value = 36
$ws.Range("D35").Value = 10066.99824741
$ws.Range("C42").Value = 30
$ws.Range("C43").Value = 11

# --- Restore the selected cell/view to match the saved state ---
$ws.Range("C26").Select()

# --- Restore the workbook window's on-screen position ---
$aw = $excel.ActiveWindow
$aw.Left = 6380
$aw.Top = 520
